$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-08 17:14:23", 0.0008),
    @("2023-12-08 17:15:47", 0.005600000000000001),
    @("2023-12-08 17:16:15", 0.002),
    @("2023-12-08 17:16:26", 0.0002)
)

$startRow = 116
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
